# Add a new daily entry row (row 45) to the "新题" sheet:
#   A45 = 2019-04-21 (date, same format as the rows above)
#   B45 = 376
#   C45 = "dp"
#   F45 = "done"
# and move the active selection to B45.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新题")

# Clone the formatting of the date cell above (A44) onto A45 first, so it
# inherits the same style (date number format) instead of Excel inventing
# a new one when we assign a date value below.
$ws.Range("A44").Copy()
$ws.Range("A45").PasteSpecial(-4122)  # xlPasteFormats

# Now fill in the values/content.
$ws.Range("A45").Value = 43576          # 2019-04-21
$ws.Range("B45").Value = 376
$ws.Range("C45").Value = "dp"
$ws.Range("F45").Value = "done"

# Match the saved selection/active cell.
$ws.Range("B45").Select()
